# Populate the previously-empty worksheet with the phone/password/name table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - write in column order so the shared-string table is
# built up in the same order as the target workbook (phone, password, name).
$ws.Range("A1").Value = "phone"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "name"

# Data row (row 2) - again write in column order (A, C, B) so any new
# shared strings land at the expected indices (لين before pass789).
$ws.Range("A2").Value = 932031600
$ws.Range("C2").Value = "لين"
$ws.Range("B2").Value = "pass789"

# Column A was best-fit to the phone-number width; emulate that sizing.
$ws.Columns.Item(1).ColumnWidth = 9.1666667

# Final UI state: active cell on A2 (matches the saved <selection>).
$ws.Range("A2").Select() | Out-Null

Write-Host "populated sheet1 with header/data rows"
